# Daily_Motivations.xlsx update: updating the daily Scores
# - Flip several existing boolean cells (0 -> 1 / False -> True) in rows 4-61
# - Append 9 new data rows (62-70) covering 2025-02-21, 2025-02-22, 2025-02-23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: correct existing boolean score cells ---
$ws.Cells(4, 14).Value = $true

$ws.Cells(7, 3).Value = $true
$ws.Cells(7, 6).Value = $true
$ws.Cells(7, 7).Value = $true
$ws.Cells(7, 9).Value = $true
$ws.Cells(7, 10).Value = $true
$ws.Cells(7, 13).Value = $true
$ws.Cells(7, 14).Value = $true

$ws.Cells(10, 3).Value = $true
$ws.Cells(10, 7).Value = $true
$ws.Cells(10, 9).Value = $true
$ws.Cells(10, 10).Value = $true
$ws.Cells(10, 12).Value = $true
$ws.Cells(10, 14).Value = $true

$ws.Cells(13, 3).Value = $true
$ws.Cells(13, 5).Value = $true
$ws.Cells(13, 7).Value = $true
$ws.Cells(13, 10).Value = $true
$ws.Cells(13, 12).Value = $true
$ws.Cells(13, 14).Value = $true

$ws.Cells(16, 7).Value = $true
$ws.Cells(16, 8).Value = $true
$ws.Cells(16, 11).Value = $true
$ws.Cells(16, 14).Value = $true

$ws.Cells(19, 8).Value = $true
$ws.Cells(19, 11).Value = $true
$ws.Cells(19, 14).Value = $true

$ws.Cells(22, 14).Value = $true

$ws.Cells(25, 14).Value = $true

$ws.Cells(28, 3).Value = $true
$ws.Cells(28, 6).Value = $true
$ws.Cells(28, 7).Value = $true
$ws.Cells(28, 9).Value = $true
$ws.Cells(28, 10).Value = $true
$ws.Cells(28, 11).Value = $true
$ws.Cells(28, 14).Value = $true

$ws.Cells(31, 3).Value = $true
$ws.Cells(31, 5).Value = $true
$ws.Cells(31, 6).Value = $true
$ws.Cells(31, 7).Value = $true
$ws.Cells(31, 9).Value = $true
$ws.Cells(31, 10).Value = $true
$ws.Cells(31, 13).Value = $true
$ws.Cells(31, 14).Value = $true

$ws.Cells(34, 3).Value = $true
$ws.Cells(34, 7).Value = $true
$ws.Cells(34, 12).Value = $true
$ws.Cells(34, 14).Value = $true

$ws.Cells(37, 3).Value = $true
$ws.Cells(37, 7).Value = $true
$ws.Cells(37, 8).Value = $true
$ws.Cells(37, 12).Value = $true
$ws.Cells(37, 14).Value = $true

$ws.Cells(40, 3).Value = $true
$ws.Cells(40, 8).Value = $true
$ws.Cells(40, 14).Value = $true

$ws.Cells(43, 14).Value = $true

$ws.Cells(46, 10).Value = $true
$ws.Cells(46, 14).Value = $true

$ws.Cells(49, 7).Value = $true
$ws.Cells(49, 10).Value = $true
$ws.Cells(49, 14).Value = $true

$ws.Cells(52, 5).Value = $true
$ws.Cells(52, 6).Value = $true
$ws.Cells(52, 7).Value = $true
$ws.Cells(52, 9).Value = $true
$ws.Cells(52, 14).Value = $true

$ws.Cells(55, 6).Value = $true
$ws.Cells(55, 7).Value = $true
$ws.Cells(55, 9).Value = $true
$ws.Cells(55, 12).Value = $true
$ws.Cells(55, 13).Value = $true
$ws.Cells(55, 14).Value = $true

$ws.Cells(58, 3).Value = $true
$ws.Cells(58, 6).Value = $true
$ws.Cells(58, 8).Value = $true
$ws.Cells(58, 12).Value = $true
$ws.Cells(58, 14).Value = $true

$ws.Cells(60, 6).Value = $true
$ws.Cells(60, 11).Value = $true

$ws.Cells(61, 3).Value = $true
$ws.Cells(61, 6).Value = $true
$ws.Cells(61, 8).Value = $true
$ws.Cells(61, 11).Value = $true
$ws.Cells(61, 12).Value = $true
$ws.Cells(61, 14).Value = $true

# --- Part 2: append new rows 62-70 ---

# Row 62
$ws.Cells(62, 1).NumberFormat = "@"
$ws.Cells(62, 1).Value = "2025-02-21"
$ws.Cells(62, 2).Value = "sleep"
$ws.Cells(62, 3).Value = $true
$ws.Cells(62, 4).Value = $false
$ws.Cells(62, 5).Value = $true
$ws.Cells(62, 6).Value = $false
$ws.Cells(62, 7).Value = $true
$ws.Cells(62, 8).Value = $true
$ws.Cells(62, 9).Value = $true
$ws.Cells(62, 10).Value = $false
$ws.Cells(62, 11).Value = $true
$ws.Cells(62, 12).Value = $true
$ws.Cells(62, 13).Value = $true
$ws.Cells(62, 14).Value = $true
$ws.Cells(62, 15).Value = $true

# Row 63
$ws.Cells(63, 1).NumberFormat = "@"
$ws.Cells(63, 1).Value = "2025-02-21"
$ws.Cells(63, 2).Value = "activity"
$ws.Cells(63, 3).Value = $true
$ws.Cells(63, 4).Value = $false
$ws.Cells(63, 5).Value = $true
$ws.Cells(63, 6).Value = $false
$ws.Cells(63, 7).Value = $true
$ws.Cells(63, 8).Value = $true
$ws.Cells(63, 9).Value = $false
$ws.Cells(63, 10).Value = $false
$ws.Cells(63, 11).Value = $false
$ws.Cells(63, 12).Value = $true
$ws.Cells(63, 13).Value = $false
$ws.Cells(63, 14).Value = $false
$ws.Cells(63, 15).Value = $false

# Row 64
$ws.Cells(64, 1).NumberFormat = "@"
$ws.Cells(64, 1).Value = "2025-02-21"
$ws.Cells(64, 2).Value = "weekly_activity"
$ws.Cells(64, 3).Value = $true
$ws.Cells(64, 4).Value = $false
$ws.Cells(64, 5).Value = $true
$ws.Cells(64, 6).Value = $true
$ws.Cells(64, 7).Value = $true
$ws.Cells(64, 8).Value = $true
$ws.Cells(64, 9).Value = $true
$ws.Cells(64, 10).Value = $true
$ws.Cells(64, 11).Value = $true
$ws.Cells(64, 12).Value = $true
$ws.Cells(64, 13).Value = $true
$ws.Cells(64, 14).Value = $true
$ws.Cells(64, 15).Value = $false

# Row 65
$ws.Cells(65, 1).NumberFormat = "@"
$ws.Cells(65, 1).Value = "2025-02-22"
$ws.Cells(65, 2).Value = "sleep"
$ws.Cells(65, 3).Value = $true
$ws.Cells(65, 4).Value = $false
$ws.Cells(65, 5).Value = $true
$ws.Cells(65, 6).Value = $true
$ws.Cells(65, 7).Value = $false
$ws.Cells(65, 8).Value = $true
$ws.Cells(65, 9).Value = $true
$ws.Cells(65, 10).Value = $true
$ws.Cells(65, 11).Value = $true
$ws.Cells(65, 12).Value = $true
$ws.Cells(65, 13).Value = $true
$ws.Cells(65, 14).Value = $true
$ws.Cells(65, 15).Value = $true

# Row 66
$ws.Cells(66, 1).NumberFormat = "@"
$ws.Cells(66, 1).Value = "2025-02-22"
$ws.Cells(66, 2).Value = "activity"
$ws.Cells(66, 3).Value = $true
$ws.Cells(66, 4).Value = $false
$ws.Cells(66, 5).Value = $false
$ws.Cells(66, 6).Value = $false
$ws.Cells(66, 7).Value = $false
$ws.Cells(66, 8).Value = $true
$ws.Cells(66, 9).Value = $true
$ws.Cells(66, 10).Value = $true
$ws.Cells(66, 11).Value = $true
$ws.Cells(66, 12).Value = $false
$ws.Cells(66, 13).Value = $false
$ws.Cells(66, 14).Value = $false
$ws.Cells(66, 15).Value = $false

# Row 67
$ws.Cells(67, 1).NumberFormat = "@"
$ws.Cells(67, 1).Value = "2025-02-22"
$ws.Cells(67, 2).Value = "weekly_activity"
$ws.Cells(67, 3).Value = $false
$ws.Cells(67, 4).Value = $false
$ws.Cells(67, 5).Value = $false
$ws.Cells(67, 6).Value = $false
$ws.Cells(67, 7).Value = $false
$ws.Cells(67, 8).Value = $false
$ws.Cells(67, 9).Value = $false
$ws.Cells(67, 10).Value = $true
$ws.Cells(67, 11).Value = $false
$ws.Cells(67, 12).Value = $false
$ws.Cells(67, 13).Value = $false
$ws.Cells(67, 14).Value = $true
$ws.Cells(67, 15).Value = $false

# Row 68
$ws.Cells(68, 1).NumberFormat = "@"
$ws.Cells(68, 1).Value = "2025-02-23"
$ws.Cells(68, 2).Value = "sleep"
$ws.Cells(68, 3).Value = $true
$ws.Cells(68, 4).Value = $false
$ws.Cells(68, 5).Value = $true
$ws.Cells(68, 6).Value = $true
$ws.Cells(68, 7).Value = $true
$ws.Cells(68, 8).Value = $true
$ws.Cells(68, 9).Value = $false
$ws.Cells(68, 10).Value = $true
$ws.Cells(68, 11).Value = $true
$ws.Cells(68, 12).Value = $true
$ws.Cells(68, 13).Value = $true
$ws.Cells(68, 14).Value = $true
$ws.Cells(68, 15).Value = $true

# Row 69
$ws.Cells(69, 1).NumberFormat = "@"
$ws.Cells(69, 1).Value = "2025-02-23"
$ws.Cells(69, 2).Value = "activity"
$ws.Cells(69, 3).Value = $false
$ws.Cells(69, 4).Value = $false
$ws.Cells(69, 5).Value = $true
$ws.Cells(69, 6).Value = $false
$ws.Cells(69, 7).Value = $true
$ws.Cells(69, 8).Value = $true
$ws.Cells(69, 9).Value = $true
$ws.Cells(69, 10).Value = $true
$ws.Cells(69, 11).Value = $false
$ws.Cells(69, 12).Value = $false
$ws.Cells(69, 13).Value = $false
$ws.Cells(69, 14).Value = $false
$ws.Cells(69, 15).Value = $false

# Row 70
$ws.Cells(70, 1).NumberFormat = "@"
$ws.Cells(70, 1).Value = "2025-02-23"
$ws.Cells(70, 2).Value = "weekly_activity"
$ws.Cells(70, 3).Value = $false
$ws.Cells(70, 4).Value = $false
$ws.Cells(70, 5).Value = $false
$ws.Cells(70, 6).Value = $false
$ws.Cells(70, 7).Value = $false
$ws.Cells(70, 8).Value = $false
$ws.Cells(70, 9).Value = $true
$ws.Cells(70, 10).Value = $true
$ws.Cells(70, 11).Value = $false
$ws.Cells(70, 12).Value = $false
$ws.Cells(70, 13).Value = $false
$ws.Cells(70, 14).Value = $true
$ws.Cells(70, 15).Value = $false

